$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the existing autofilter criteria (column A == "Sow") which unhides
# all the previously filtered-out rows and clears FilterMode.
$ws.AutoFilterMode = $false

# Append the new "Tags" row (row 89) with its parameter metadata.
$ws.Range("A89").Value2 = "Tags"
$ws.Range("B89").Value2 = "extra_tops"
$ws.Range("C89").Value2 = "Additional Top Level Tags"
$ws.Range("D89").Value2 = "_"
$ws.Range("E89").Value2 = 300
$ws.Range("F89").Value2 = "TextDict"
$ws.Range("H89").Value2 = 0
$ws.Range("I89").Value2 = 0

# Re-apply the autofilter over the full (now larger) range, with no active
# filter criteria, so every row is shown.
$ws.Range("A1:I89").AutoFilter()

# Keep the hidden _FilterDatabase defined name in sync with the new range.
for ($i = 1; $i -le $wb.Names.Count; $i++) {
  $n = $wb.Names.Item($i)
  if ($n.Name -like "*_FilterDatabase*") {
    $n.RefersTo = "=game_params!`$A`$1:`$I`$89"
  }
}

# Update the active selection to match the edited cell.
$ws.Range("C89").Select()
